$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

# Delete row 77 (Caso -657, "Conde 1632") - this shifts subsequent rows up by one
$ws.Rows.Item(77).Delete()
